# Applies "more updates to get sdusd sim to run" edits to Sheet1:
# - Q2 value changes from 3 to 5
# - Three new simulation rows (3,4,5) + a fourth variant (row 6) are appended,
#   varying F (min_diff), G (max_diff) and, for the last row, O (ta_sd).
# - Selection/view moves to scroll right and land on Q7.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Fix existing row 2: n_cohorts (Q2) 3 -> 5 ---
$ws.Range("Q2").Value = 5

# --- Template values shared by the new rows (same as row 2 except where noted) ---
# A..Z -> nsims, npoints, test_SEM, impact_type, impact_function, min_diff, max_diff,
#         weight_type, method, weighted_average, covariates, peer_effects, stud_sorting,
#         rho, ta_sd, tc_sd, n_cohorts, pretest_coef, num_cats, lin_alpha, pctile,
#         weight_below, weight_above, v_alpha, mrpctile, mrdist

$rowsData = @(
    # row, F(min_diff), G(max_diff), O(ta_sd)
    @{ Row = 3; F = 0;   G = 0.3; O = 0.6 },
    @{ Row = 4; F = 0.1; G = 0.3; O = 0.6 },
    @{ Row = 5; F = 0.2; G = 0.4; O = 0.6 },
    @{ Row = 6; F = 0.2; G = 0.4; O = 0.1 }
)

foreach ($rd in $rowsData) {
    $r = $rd.Row

    $ws.Cells.Item($r, 1).Value = 20
    $ws.Cells.Item($r, 2).Value = 500
    $ws.Cells.Item($r, 3).Value = 0.07
    $ws.Cells.Item($r, 4).Value = "No"
    $ws.Cells.Item($r, 5).Value = 1
    $ws.Cells.Item($r, 6).Value = $rd.F
    $ws.Cells.Item($r, 7).Value = $rd.G
    $ws.Cells.Item($r, 8).Value = "rawlsian"
    $ws.Cells.Item($r, 9).Value = "bin"
    $ws.Cells.Item($r, 10).Value = 0
    $ws.Cells.Item($r, 11).Value = 0
    $ws.Cells.Item($r, 12).Value = 0
    $ws.Cells.Item($r, 13).Value = 0
    $ws.Cells.Item($r, 14).Value = 0.3
    $ws.Cells.Item($r, 15).Value = $rd.O
    $ws.Cells.Item($r, 16).Value = 1
    $ws.Cells.Item($r, 17).Value = 5
    $ws.Cells.Item($r, 18).Value = 0.9
    $ws.Cells.Item($r, 19).Value = 2
    $ws.Cells.Item($r, 20).Value = 2
    $ws.Cells.Item($r, 21).Value = 0.4
    $ws.Cells.Item($r, 22).Value = 0.8
    $ws.Cells.Item($r, 23).Value = 0.2
    $ws.Cells.Item($r, 24).Value = 1
    $ws.Cells.Item($r, 25).Value = 0.4
    $ws.Cells.Item($r, 26).Value = 0.2
}

# --- Update the view: scroll right to column M and select Q7 ---
$ws.Activate()
try { $excel.ActiveWindow.ScrollColumn = 13 } catch {}
try { $excel.ActiveWindow.ScrollRow = 1 } catch {}
$ws.Range("Q7").Select()
